$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K to F:M)
$ws.Range("D:E").Insert()

# Copy number formatting/style from column F (the shifted-original column D) into the new D:E columns
# (restricted to the rows that actually carry data, so label-only rows 5/6/37/79 stay untouched)
$ws.Range("F7:F35,F38:F77,F80:F102").Copy()
$ws.Range("D7:E35,D38:E77,D80:E102").PasteSpecial(-4122)

# Fill the two new columns with the new quarter data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 910100
$ws.Range("E8").Value = 881300
$ws.Range("D9").Value = 655600
$ws.Range("E9").Value = 628000
$ws.Range("D10").Value = 254500
$ws.Range("E10").Value = 253300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 1800
$ws.Range("E14").Value = 2600
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 805700
$ws.Range("E17").Value = 772600
$ws.Range("D18").Value = 104400
$ws.Range("E18").Value = 108700
$ws.Range("D20").Value = -6000
$ws.Range("E20").Value = 1100
$ws.Range("D21").Value = 138500
$ws.Range("E21").Value = 144900
$ws.Range("D22").Value = 18500
$ws.Range("E22").Value = 12500
$ws.Range("D23").Value = 79900
$ws.Range("E23").Value = 97300
$ws.Range("D24").Value = 19100
$ws.Range("E24").Value = 25000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 60800
$ws.Range("E26").Value = 72300
$ws.Range("D27").Value = 60000
$ws.Range("E27").Value = 71600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 6000
$ws.Range("E32").Value = -1100
$ws.Range("D33").Value = 60000
$ws.Range("E33").Value = 71600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 60000
$ws.Range("E35").Value = 71600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 132500
$ws.Range("E41").Value = 153700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 663200
$ws.Range("E43").Value = 685900
$ws.Range("D44").Value = 835700
$ws.Range("E44").Value = 841000
$ws.Range("D45").Value = 105800
$ws.Range("E45").Value = 105700
$ws.Range("D46").Value = 1737200
$ws.Range("E46").Value = 1786300
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("D48").Value = 912100
$ws.Range("E48").Value = 886800
$ws.Range("D49").Value = 1693700
$ws.Range("E49").Value = 1708500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 102200
$ws.Range("E52").Value = 117900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4445200
$ws.Range("E54").Value = 4499500
$ws.Range("D57").Value = 273200
$ws.Range("E57").Value = 282800
$ws.Range("D58").Value = 43000
$ws.Range("E58").Value = 48400
$ws.Range("D59").Value = 369400
$ws.Range("E59").Value = 325400
$ws.Range("D60").Value = 685600
$ws.Range("E60").Value = 656600
$ws.Range("D61").Value = 1638600
$ws.Range("E61").Value = 1681700
$ws.Range("D62").Value = 478300
$ws.Range("E62").Value = 520800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2865600
$ws.Range("E66").Value = 2919600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1630200
$ws.Range("E72").Value = 1595400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1579600
$ws.Range("E76").Value = 1579900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 60000
$ws.Range("E81").Value = 71600
$ws.Range("D83").Value = 40100
$ws.Range("E83").Value = 35100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 137500
$ws.Range("E89").Value = 137200
$ws.Range("D91").Value = -49800
$ws.Range("E91").Value = -23200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -54900
$ws.Range("E94").Value = -774300
$ws.Range("D96").Value = -21500
$ws.Range("E96").Value = -21500
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -104200
$ws.Range("E100").Value = 649400
$ws.Range("D101").Value = -300
$ws.Range("E101").Value = -3900
$ws.Range("D102").Value = -21900
$ws.Range("E102").Value = 8400
